$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data scraped on Mon Sep 11 18:27:05 UTC 2023

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.286.48'
$ws.Range('E2').Value = '  -2.57%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.566.71'
$ws.Range('E3').Value = '  -3.59%  '

$ws.Range('E4').Value = '  -0.34%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '207.51'
$ws.Range('E5').Value = '  -2.96%  '

$ws.Range('E6').Value = '  -0.33%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.478'
$ws.Range('E7').Value = '  -5.06%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.0608'
$ws.Range('E8').Value = '  -1.60%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.242'
$ws.Range('E9').Value = '  -2.89%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '17.80'
$ws.Range('E10').Value = '  -2.66%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0782'
$ws.Range('E11').Value = '  -0.67%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.784.95'
$ws.Range('E12').Value = '  -3.52%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.585.62'
$ws.Range('E13').Value = '  -2.42%  '

$ws.Range('E14').Value = '  -3.88%  '

$ws.Range('E15').Value = '  -3.47%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '25.291.74'
$ws.Range('E16').Value = '  -2.45%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '59.47'
$ws.Range('E17').Value = '  -2.91%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0₃0713'
$ws.Range('E18').Value = '  -3.26%  '

$ws.Range('E19').Value = '  -0.32%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '185.63'
$ws.Range('E20').Value = '  -3.01%  '

$ws.Range('E21').Value = '  -2.54%  '

$ws.Range('E22').Value = '  -3.03%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.88'
$ws.Range('E23').Value = '  -3.18%  '

$ws.Range('B24').Value = 'BinanceUSD'
$ws.Range('C24').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.00'
$ws.Range('E24').Value = '  -0.33%  '

$ws.Range('B25').Value = 'Stellar'
$ws.Range('C25').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.128'
$ws.Range('E25').Value = '  -4.17%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '140.07'
$ws.Range('E26').Value = '  -2.41%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.66'
$ws.Range('E27').Value = '  -5.78%  '

$ws.Range('E28').Value = '  -3.93%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '14.85'
$ws.Range('E29').Value = '  -2.23%  '

$ws.Range('E30').Value = '  -5.93%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0463'
$ws.Range('E31').Value = '  -3.99%  '

$ws.Range('E32').Value = '  -3.06%  '

$ws.Range('E33').Value = '  -3.44%  '

$ws.Range('E34').Value = '  -2.09%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.31'
$ws.Range('E35').Value = '  -3.69%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.088.46'
$ws.Range('E36').Value = '  -3.22%  '

$ws.Range('E37').Value = '  -0.62%  '

$ws.Range('E38').Value = '  -4.99%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0149'
$ws.Range('E39').Value = '  -2.75%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.495'
$ws.Range('E40').Value = '  -4.39%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.771'
$ws.Range('E41').Value = '  -8.89%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.772'
$ws.Range('E42').Value = '  +0.72%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '93.32'
$ws.Range('E43').Value = '  -4.59%  '

$ws.Range('E44').Value = '  -1.94%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.699.22'
$ws.Range('E45').Value = '  -3.48%  '

$ws.Range('E46').Value = '  -2.60%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '52.72'
$ws.Range('E47').Value = '  -3.30%  '

$ws.Range('E48').Value = '  -4.61%  '

$ws.Range('E49').Value = '  -2.48%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.406'
$ws.Range('E50').Value = '  -1.66%  '

$ws.Range('E51').Value = '  -0.53%  '
